# Updated symbol list (crypto price/volume refresh) applied via COM automation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell reference -> new value (Price / Volume(1h) columns).
# Every one of these cells already holds its number formatted as literal
# text (e.g. "0.1540", "277.09", "0.94%") rather than a real Number, so we
# force the Text number format before writing the replacement string. That
# keeps trailing zeros (e.g. "0.08980") and the literal "%" characters intact
# instead of letting Excel re-interpret the input as a Number/Percentage.
$updates = [ordered]@{
    "D2" = "276.83"
    "E2" = "0.79%"
    "D3" = "27.34"
    "E3" = "1.89%"
    "D4" = "4.873"
    "E4" = "-0.72%"
    "D5" = "0.06415"
    "E5" = "1.35%"
    "D6" = "6.942"
    "E6" = "1.24%"
    "D7" = "1.182"
    "E7" = "-5.73%"
    "D8" = "0.8761"
    "E8" = "0.77%"
    "D9" = "0.1534"
    "E9" = "-11.37%"
    "D10" = "0.05160"
    "E10" = "2.67%"
    "D11" = "0.07474"
    "E11" = "0.69%"
    "D12" = "0.02951"
    "E12" = "-0.58%"
    "D13" = "0.08980"
    "E13" = "-0.42%"
    "D14" = "0.001570"
    "E14" = "-0.04%"
    "D15" = "0.0006360"
    "E15" = "0.91%"
    "D16" = "0.006124"
    "E16" = "4.31%"
    "D17" = "3.481"
    "E17" = "1.01%"
    "E18" = "-0.17%"
    "E19" = "0.11%"
    "E20" = "-0.08%"
    "E21" = "-1.06%"
    "D22" = "3.907"
    "E22" = "-0.02%"
    "D23" = "0.04402"
    "E23" = "1.19%"
    "D24" = "0.1499"
    "E24" = "8.62%"
    "D26" = "0.001177"
    "E26" = "-0.01%"
    "E27" = "-8.97%"
    "E28" = "8.27%"
    "E29" = "15.01%"
    "D40" = "0.04165"
    "E40" = "3.05%"
    "D41" = "0.006802"
    "E41" = "1.75%"
    "E42" = "0.63%"
    "D43" = "0.002040"
    "E43" = "-6.87%"
    "E44" = "6.86%"
    "D45" = "0.00005302"
    "E45" = "-0.01%"
    "D46" = "1.685"
    "E46" = "13.10%"
    "D47" = "0.01851"
    "E47" = "-11.92%"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}

